$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric data values (columns B-G, rows 2-12)
$data = @{
    2  = @(0, 0.061133388136613197, 0, 0.061805075695871202, 0, 0.065184352889597494)
    3  = @(0.0071677931762149398, 0.062207675899909998, 0.00709254656995885, 0.062722743604473705, 0.0067998303527872098, 0.065452920276735202)
    4  = @(0.0143355863469722, 0.063267528295928896, 0.0141850931435136, 0.063640411512031703, 0.013599660702574101, 0.065691102055548997)
    5  = @(0.021503379525753801, 0.064320593712562199, 0.021277639715795899, 0.064558079425248396, 0.020399491052263199, 0.065929283828159296)
    6  = @(0.028671172693532001, 0.065373659137109297, 0.0283701862895827, 0.065468858483853601, 0.027199321358782801, 0.066167465607028506)
    7  = @(0.035838965867301897, 0.066421700353464896, 0.035462732857769998, 0.066378585567892806, 0.033999151758236201, 0.066405647385897604)
    8  = @(0.043006759040784798, 0.067450747510531706, 0.042555279431814703, 0.067270312831190601, 0.0407989821063783, 0.066643829164766799)
    9  = @(0.050174552212917098, 0.068422117376008104, 0.049647826004628401, 0.068115788669579402, 0.047598812423613898, 0.066882010943630901)
    10 = @(0.057342345387291098, 0.069393487242292301, 0.056740372566541498, 0.068961264507401199, 0.054398642815280902, 0.067120192722247604)
    11 = @(0.064510138558826299, 0.0703648571161644, 0.063832919147396597, 0.069806740343181697, 0.061198473162006002, 0.067358374507713106)
    12 = @(0.071677931734810199, 0.071221169955120706, 0.070925465719669206, 0.0704231366345692, 0.067998303518177705, 0.067573826441724494)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # B=2 .. G=7
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# Set explicit column widths (B-G) to approximate the bestFit results from the diff
# (engine quantizes ColumnWidth to 1/6 units; these inputs round-trip to the
# closest achievable widths to the authored bestFit values)
$ws.Columns.Item(2).ColumnWidth = 12.333333333333334
$ws.Columns.Item(3).ColumnWidth = 11.0
$ws.Columns.Item(4).ColumnWidth = 11.333333333333334
$ws.Columns.Item(5).ColumnWidth = 11.0
$ws.Columns.Item(6).ColumnWidth = 11.666666666666666
$ws.Columns.Item(7).ColumnWidth = 11.0

$wb.Save()
